$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 463, shifting existing rows 463-482 down to 464-483
$ws.Rows(463).Insert()

$ws.Range("A463").Value = 10
$ws.Range("B463").Value = "Vega Modelo de Temuco"
$ws.Range("C463").Value = "La Araucanía"
$ws.Range("D463").Value = 45147
$ws.Range("D463").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E463").Value = 9
$ws.Range("F463").Value = "Fruta"
$ws.Range("G463").Value = 100102
$ws.Range("H463").Value = "Cítricos"
$ws.Range("I463").Value = 100102006
$ws.Range("J463").Value = "Pomelo"
$ws.Range("K463").Value = "Start Ruby"
$ws.Range("L463").Value = "Primera"
$ws.Range("M463").Value = 90
$ws.Range("N463").Value = 14000
$ws.Range("O463").Value = 15000
$ws.Range("P463").Value = 14556
$ws.Range("Q463").Value = "$/bandeja 15 kilos granel"
$ws.Range("R463").Value = "Región de O'Higgins"
$ws.Range("S463").Value = 970
$ws.Range("T463").Value = 15
